$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (110 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 753.6177
$ws.Range("I15").Value = 753.6177
$ws.Range("K15").Value = 2260.8531
$ws.Range("M15").Value = -2091.8531
$ws.Range("H28").Value = 3517.4375
$ws.Range("I28").Value = 948.875
$ws.Range("J28").Value = 6086
$ws.Range("K28").Value = 948.875
$ws.Range("L28").Value = 6086
$ws.Range("M28").Value = -463.875
$ws.Range("N28").Value = -7056
$ws.Range("H38").Value = 3412.6667
$ws.Range("I38").Value = 2216.889
$ws.Range("K38").Value = 6650.667
$ws.Range("M38").Value = -6278.667
$ws.Range("H39").Value = 207.78572
$ws.Range("I39").Value = 113.22222
$ws.Range("J39").Value = 378
$ws.Range("K39").Value = 339.66666
$ws.Range("L39").Value = 1134
$ws.Range("M39").Value = -43.66665999999998
$ws.Range("N39").Value = -1726
$ws.Range("H48").Value = 6315
$ws.Range("J48").Value = 6315
$ws.Range("L48").Value = 18945
$ws.Range("N48").Value = -19529
$ws.Range("H56").Value = 6315
$ws.Range("J56").Value = 6315
$ws.Range("L56").Value = 18945
$ws.Range("N56").Value = -20013
$ws.Range("H62").Value = 8128
$ws.Range("J62").Value = 9374.25
$ws.Range("L62").Value = 9374.25
$ws.Range("N62").Value = -10622.25
$ws.Range("H65").Value = 8128
$ws.Range("J65").Value = 9374.25
$ws.Range("L65").Value = 46871.25
$ws.Range("N65").Value = -53111.25
$ws.Range("H86").Value = 4228.4287
$ws.Range("I86").Value = 3919.8
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3919.8
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2796.8
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4228.4287
$ws.Range("I89").Value = 3919.8
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 19599
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -13983
$ws.Range("N89").Value = -36232
$ws.Range("H92").Value = 144.36363
$ws.Range("I92").Value = 108.8
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 108.8
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 1139.2
$ws.Range("N92").Value = -2996
$ws.Range("H98").Value = 528.13336
$ws.Range("I98").Value = 494.42856
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 494.42856
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 1003.57144
$ws.Range("N98").Value = -3996
$ws.Range("H107").Value = 817.94116
$ws.Range("I107").Value = 940.3333
$ws.Range("J107").Value = 524.2
$ws.Range("K107").Value = 940.3333
$ws.Range("L107").Value = 524.2
$ws.Range("M107").Value = 979.6667
$ws.Range("N107").Value = -4364.2
$ws.Range("H112").Value = 2077
$ws.Range("J112").Value = 1971.25
$ws.Range("L112").Value = 5913.75
$ws.Range("N112").Value = -8129.75
$ws.Range("H113").Value = 3066.3333
$ws.Range("I113").Value = 2849.75
$ws.Range("K113").Value = 2849.75
$ws.Range("M113").Value = 404.25
$ws.Range("H121").Value = 700
$ws.Range("J121").Value = 700
$ws.Range("L121").Value = 2100
$ws.Range("N121").Value = -5594
$ws.Range("H122").Value = 528.13336
$ws.Range("I122").Value = 494.42856
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 1483.28568
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 966.71432
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 19039.5
$ws.Range("I132").Value = 19039.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 57118.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -54588.5
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 4245.125
$ws.Range("I138").Value = 3565.6667
$ws.Range("J138").Value = 4652.8
$ws.Range("K138").Value = 10697.0001
$ws.Range("L138").Value = 13958.4
$ws.Range("M138").Value = -5557.000100000001
$ws.Range("N138").Value = -24238.4
$ws.Range("H141").Value = 1631.9166
$ws.Range("I141").Value = 1058.3
$ws.Range("K141").Value = 3174.9
$ws.Range("M141").Value = 2005.1

# --- Sheet: ARM (48 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1229.8889
$ws.Range("I2").Value = 1008.75
$ws.Range("K2").Value = 1008.75
$ws.Range("M2").Value = -895.75
$ws.Range("H13").Value = 1090.6
$ws.Range("I13").Value = 1151
$ws.Range("K13").Value = 1151
$ws.Range("M13").Value = -1007
$ws.Range("H45").Value = 3055.2856
$ws.Range("I45").Value = 2666.3333
$ws.Range("K45").Value = 2666.3333
$ws.Range("M45").Value = -2289.3333
$ws.Range("H61").Value = 3338.9333
$ws.Range("I61").Value = 2698.7693
$ws.Range("K61").Value = 2698.7693
$ws.Range("M61").Value = -2486.7693
$ws.Range("H74").Value = 4998.25
$ws.Range("J74").Value = 3700
$ws.Range("L74").Value = 3700
$ws.Range("N74").Value = -5448
$ws.Range("H77").Value = 4998.25
$ws.Range("J77").Value = 3700
$ws.Range("L77").Value = 18500
$ws.Range("N77").Value = -27236
$ws.Range("H110").Value = 846.1539
$ws.Range("I110").Value = 850.1
$ws.Range("K110").Value = 850.1
$ws.Range("M110").Value = 1194.9
$ws.Range("H116").Value = 1229.8889
$ws.Range("I116").Value = 1008.75
$ws.Range("K116").Value = 1008.75
$ws.Range("M116").Value = 1285.25
$ws.Range("H122").Value = 1903.5
$ws.Range("I122").Value = 1344.2
$ws.Range("K122").Value = 4032.6
$ws.Range("M122").Value = -1582.6
$ws.Range("H132").Value = 5329.8
$ws.Range("I132").Value = 4858.647
$ws.Range("K132").Value = 14575.941
$ws.Range("M132").Value = -12045.941
$ws.Range("H135").Value = 26952.334
$ws.Range("J135").Value = 22929
$ws.Range("L135").Value = 22929
$ws.Range("N135").Value = -33069
$ws.Range("H136").Value = 3338.9333
$ws.Range("I136").Value = 2698.7693
$ws.Range("K136").Value = 8096.3079
$ws.Range("M136").Value = -5546.3079

# --- Sheet: BSM (47 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1229.8889
$ws.Range("I3").Value = 1008.75
$ws.Range("K3").Value = 1008.75
$ws.Range("M3").Value = -894.75
$ws.Range("H26").Value = 16805.25
$ws.Range("I26").Value = 16805.25
$ws.Range("K26").Value = 16805.25
$ws.Range("M26").Value = -16513.25
$ws.Range("H80").Value = 986.0833
$ws.Range("I80").Value = 767.5714
$ws.Range("J80").Value = 1292
$ws.Range("K80").Value = 767.5714
$ws.Range("L80").Value = 1292
$ws.Range("M80").Value = 230.4286
$ws.Range("N80").Value = -3288
$ws.Range("H83").Value = 986.0833
$ws.Range("I83").Value = 767.5714
$ws.Range("J83").Value = 1292
$ws.Range("K83").Value = 3837.857
$ws.Range("L83").Value = 6460
$ws.Range("M83").Value = 1154.143
$ws.Range("N83").Value = -16444
$ws.Range("H99").Value = 2072.1428
$ws.Range("I99").Value = 1877.5
$ws.Range("J99").Value = 2331.6667
$ws.Range("K99").Value = 1877.5
$ws.Range("L99").Value = 2331.6667
$ws.Range("M99").Value = -379.5
$ws.Range("N99").Value = -5327.6667
$ws.Range("H105").Value = 2939
$ws.Range("I105").Value = 2600.1667
$ws.Range("J105").Value = 3955.5
$ws.Range("K105").Value = 2600.1667
$ws.Range("L105").Value = 3955.5
$ws.Range("M105").Value = -853.1667000000002
$ws.Range("N105").Value = -7449.5
$ws.Range("H107").Value = 4174.2
$ws.Range("I107").Value = 3041.4
$ws.Range("J107").Value = 6439.8
$ws.Range("K107").Value = 3041.4
$ws.Range("L107").Value = 6439.8
$ws.Range("M107").Value = -1121.4
$ws.Range("N107").Value = -10279.8
$ws.Range("H134").Value = 2750
$ws.Range("I134").Value = 2750
$ws.Range("K134").Value = 8250
$ws.Range("M134").Value = -5715

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 580.38464
$ws.Range("I16").Value = 505.55554
$ws.Range("K16").Value = 505.55554
$ws.Range("M16").Value = -218.55554
$ws.Range("H113").Value = 580.38464
$ws.Range("I113").Value = 505.55554
$ws.Range("K113").Value = 505.55554
$ws.Range("M113").Value = 1664.44446
$ws.Range("H117").Value = 25356
$ws.Range("I117").Value = 20000
$ws.Range("J117").Value = 30712
$ws.Range("K117").Value = 20000
$ws.Range("L117").Value = 30712
$ws.Range("M117").Value = -15411
$ws.Range("N117").Value = -39890
$ws.Range("H132").Value = 4359.6665
$ws.Range("I132").Value = 2090.1428
$ws.Range("J132").Value = 7537
$ws.Range("K132").Value = 6270.428400000001
$ws.Range("L132").Value = 22611
$ws.Range("M132").Value = -3740.428400000001
$ws.Range("N132").Value = -27671
$ws.Range("H134").Value = 1974.8572
$ws.Range("I134").Value = 2745
$ws.Range("J134").Value = 1397.25
$ws.Range("K134").Value = 8235
$ws.Range("L134").Value = 4191.75
$ws.Range("M134").Value = -5700
$ws.Range("N134").Value = -9261.75

# --- Sheet: CUL (32 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4473.7856
$ws.Range("I17").Value = 73.5
$ws.Range("J17").Value = 6233.9
$ws.Range("K17").Value = 220.5
$ws.Range("L17").Value = 18701.7
$ws.Range("M17").Value = -51.5
$ws.Range("N17").Value = -19039.7
$ws.Range("H38").Value = 134.28572
$ws.Range("J38").Value = 268.33334
$ws.Range("L38").Value = 805.0000200000001
$ws.Range("N38").Value = -1499.00002
$ws.Range("H86").Value = 453.25
$ws.Range("I86").Value = 350
$ws.Range("J86").Value = 487.66666
$ws.Range("K86").Value = 1050
$ws.Range("L86").Value = 1462.99998
$ws.Range("M86").Value = 136
$ws.Range("N86").Value = -3834.99998
$ws.Range("H89").Value = 453.25
$ws.Range("I89").Value = 350
$ws.Range("J89").Value = 487.66666
$ws.Range("K89").Value = 3150
$ws.Range("L89").Value = 4388.99994
$ws.Range("M89").Value = 2778
$ws.Range("N89").Value = -16244.99994
$ws.Range("H112").Value = 1561.3334
$ws.Range("I112").Value = 1227
$ws.Range("J112").Value = 2230
$ws.Range("K112").Value = 3681
$ws.Range("L112").Value = 6690
$ws.Range("M112").Value = -2573
$ws.Range("N112").Value = -8906

# --- Sheet: GSM (33 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13224444
$ws.Range("I11").Value = 12860000
$ws.Range("J11").Value = 14500000
$ws.Range("K11").Value = 12860000
$ws.Range("L11").Value = 14500000
$ws.Range("M11").Value = -12859861
$ws.Range("N11").Value = -14500278
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -1707
$ws.Range("N18").Value = -2586
$ws.Range("H55").Value = 29630
$ws.Range("I55").Value = 29630
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 29630
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -29303
$ws.Range("N55").ClearContents()
$ws.Range("H102").Value = 2787.625
$ws.Range("I102").Value = 2558.5
$ws.Range("K102").Value = 2558.5
$ws.Range("M102").Value = -936.5
$ws.Range("H122").Value = 3332
$ws.Range("I122").Value = 3332
$ws.Range("K122").Value = 9996
$ws.Range("M122").Value = -7546
$ws.Range("H132").Value = 3029.9285
$ws.Range("I132").Value = 2570.7693
$ws.Range("K132").Value = 7712.3079
$ws.Range("M132").Value = -5182.3079

# --- Sheet: LTW (53 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 701
$ws.Range("J12").Value = 701
$ws.Range("L12").Value = 701
$ws.Range("N12").Value = -1041
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 3845
$ws.Range("I40").Value = 3845
$ws.Range("K40").Value = 3845
$ws.Range("M40").Value = -3709
$ws.Range("H61").Value = 3389.1667
$ws.Range("I61").Value = 1808.75
$ws.Range("K61").Value = 1808.75
$ws.Range("M61").Value = -1606.75
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 8862.375
$ws.Range("I100").Value = 5449.5
$ws.Range("K100").Value = 5449.5
$ws.Range("M100").Value = -4908.5
$ws.Range("H113").Value = 3389.1667
$ws.Range("I113").Value = 1808.75
$ws.Range("K113").Value = 1808.75
$ws.Range("M113").Value = 361.25
$ws.Range("H122").Value = 3960.75
$ws.Range("I122").Value = 3960.75
$ws.Range("K122").Value = 11882.25
$ws.Range("M122").Value = -9432.25
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet: WVR (20 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 749.5
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 500
$ws.Range("M29").Value = -210
$ws.Range("H62").Value = 11500
$ws.Range("J62").Value = 11500
$ws.Range("L62").Value = 11500
$ws.Range("N62").Value = -12748
$ws.Range("H65").Value = 11500
$ws.Range("J65").Value = 11500
$ws.Range("L65").Value = 57500
$ws.Range("N65").Value = -63740
$ws.Range("H113").Value = 885.61536
$ws.Range("I113").Value = 734
$ws.Range("K113").Value = 2202
$ws.Range("M113").Value = -32
$ws.Range("H132").Value = 2905.6667
$ws.Range("I132").Value = 2897.0908
$ws.Range("K132").Value = 8691.2724
$ws.Range("M132").Value = -6161.2724
